# The site rebuild dropped the trailing "Ver no Jupiter ..." / copyright
# footer paragraphs (along with the blank paragraph that preceded them)
# from the end of the document, right after the "Requisitos" section's
# "LOB1004: Cálculo II (Requisito fraco)" line.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph by its text.
$jupiterPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $jupiterPara = $p
        break
    }
}

if ($jupiterPara -eq $null) {
    Write-Output "Target paragraph not found; no changes made."
} else {
    # The blank paragraph right before it, and the copyright paragraph
    # right after it, are removed together with it.
    $blankPara = $jupiterPara.Previous()
    $copyrightPara = $jupiterPara.Next()

    $startPos = $blankPara.Range.Start
    $endPos = $copyrightPara.Range.End

    $r = $d.Range($startPos, $endPos)
    $r.Delete()

    Write-Output "Removed footer paragraphs."
}
